# Update row 2 of the review database sheet with new appid/keyword/review values,
# then leave the active selection on A2 (matching the saved sheetView state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "block.chain.technology"
$ws.Range("B2").Value = "blockchain technology"
$ws.Range("F2").Value = "very good app"

$ws.Range("A2").Select()
